# ---------------------------------------------------------------------------
# Add the "2022-Q4" sheet (new quarter of fund-holdings data) positioned
# right after the "总计" (total) sheet and before "2022-Q3", then backfill
# the new-sheet data and update the "总计" summary sheet with a matching
# summary row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Insert the new worksheet in the correct tab position -------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")
$q4Sheet    = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# ---- 2. Header row (bold, centered, thin border) --------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $hc = $q4Sheet.Cells.Item(1, $i + 2)
    $hc.Value = $headers[$i]
    $hc.Font.Bold = $true
    $hc.HorizontalAlignment = -4108
    $hc.VerticalAlignment = -4160
    for ($edge = 1; $edge -le 4; $edge++) {
        $hc.Borders.Item($edge).LineStyle = 1
    }
}

# ---- 3. Fund-holdings data rows (A=index, B..G mostly text, H=rank) ------
$q4rows = @(
    @(0, '012930', '中庚价值先锋股票', '68.71', '94.78', '4.39', '3.0164', 7),
    @(1, '160106', '南方高增长混合（LOF）', '16.27', '88.50', '4.92', '0.8005', 7),
    @(2, '014029', '浦银安盛红利精选混合C', '9.07', '71.06', '4.79', '0.4345', 2),
    @(3, '519115', '浦银安盛红利精选混合A', '7.60', '71.06', '4.79', '0.3640', 2),
    @(4, '010852', '中欧内需成长混合A', '3.73', '93.04', '8.73', '0.3256', 2),
    @(5, '160105', '南方积极配置混合（LOF）', '5.52', '89.92', '5.14', '0.2837', 6),
    @(6, '450010', '国富策略回报混合', '9.53', '79.59', '1.84', '0.1754', 8),
    @(7, '519125', '浦银安盛消费升级混合A', '2.11', '82.27', '8.02', '0.1692', 2),
    @(8, '005620', '中欧品质消费股票A', '1.65', '93.90', '9.01', '0.1487', 1),
    @(9, '005621', '中欧品质消费股票C', '1.01', '93.90', '9.01', '0.0910', 1),
    @(10, '000554', '南方中国梦灵活配置混合', '1.50', '92.96', '5.43', '0.0814', 7),
    @(11, '519176', '浦银安盛消费升级混合C', '0.92', '82.27', '8.02', '0.0738', 2),
    @(12, '000974', '安信消费医药主题股票', '2.33', '92.45', '2.85', '0.0664', 7),
    @(13, '010853', '中欧内需成长混合C', '0.54', '93.04', '8.73', '0.0471', 2),
    @(14, '009206', '兴银丰运稳益回报混合C', '3.03', '39.08', '1.12', '0.0339', 9),
    @(15, '011858', '安信消费升级一年持有期混合A', '0.84', '88.78', '3.88', '0.0326', 3),
    @(16, '009205', '兴银丰运稳益回报混合A', '1.91', '39.08', '1.12', '0.0214', 9),
    @(17, '011859', '安信消费升级一年持有期混合C', '0.13', '88.78', '3.88', '0.0050', 3),
    @(18, '000761', '国富健康优质生活股票', '0.11', '86.45', '3.55', '0.0039', 6),
    @(19, '011771', '国寿安保稳隆混合A', '0.50', '32.85', '0.76', '0.0038', 9),
    @(20, '001932', '国寿安保灵活优选混合', '0.11', '39.50', '1.01', '0.0011', 10),
    @(21, '011772', '国寿安保稳隆混合C', '0.00', '32.85', '0.76', 0, 9)
)

$r = 2
foreach ($row in $q4rows) {
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]

    $cCode = $q4Sheet.Cells.Item($r, 2)
    $cCode.NumberFormat = "@"
    $cCode.Value = $row[1]

    $q4Sheet.Cells.Item($r, 3).Value = $row[2]

    $cScale = $q4Sheet.Cells.Item($r, 4)
    $cScale.NumberFormat = "@"
    $cScale.Value = $row[3]

    $cPos = $q4Sheet.Cells.Item($r, 5)
    $cPos.NumberFormat = "@"
    $cPos.Value = $row[4]

    $cPct = $q4Sheet.Cells.Item($r, 6)
    $cPct.NumberFormat = "@"
    $cPct.Value = $row[5]

    $cVal = $q4Sheet.Cells.Item($r, 7)
    if ($row[6] -eq 0) {
        $cVal.Value = $row[6]
    } else {
        $cVal.NumberFormat = "@"
        $cVal.Value = $row[6]
    }

    $q4Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---- 4. Update the "总计" (total) summary sheet ---------------------------
# Insert a new row 2 for the 2022-Q4 summary figures, then renumber the
# existing A-column sequence index (0-based) for the rows that shifted down.
$totalSheet.Rows("2:2").Insert()
# The freshly inserted row inherits stray formatting from the row above it;
# strip that back to the sheet default before re-applying the one style the
# "A" column is actually supposed to carry (bold/centered/bordered, same as
# every other row-index cell in this column).
$totalSheet.Range("A2:D2").ClearFormats()

$aCell = $totalSheet.Cells.Item(2, 1)
$aCell.Value = 0
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160
for ($edge = 1; $edge -le 4; $edge++) {
    $aCell.Borders.Item($edge).LineStyle = 1
}

$cQ4 = $totalSheet.Cells.Item(2, 2)
$cQ4.NumberFormat = "@"
$cQ4.Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 22
$totalSheet.Cells.Item(2, 4).Value = 6.18

for ($row = 3; $row -le 8; $row++) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
}

# ---- 5. Keep the original active-sheet selection (总计) ------------------
$totalSheet.Activate()

Write-Host "2022-Q4 sheet added; 总计 sheet updated."
